$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data refresh: previous "Estado de Cuenta" (EC) periods are replaced
# --- The period list (column E, rows 16-21) is now listed in descending
# --- order (1906 .. 1901) and the "Valor Mora" (column F) values travel
# --- with their respective period.
$periods = @("1906", "1905", "1904", "1903", "1902", "1901")
$valores = @(17708, 31249, 31249, 31249, 31249, 31249)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}

# --- Column widths were refreshed (best-fit) after the data update.
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
